$d = $word.ActiveDocument

function Set-ParaFont($para, $bold) {
    $para.Range.Font.Name = "Helvetica"
    $para.Range.Font.Size = 10
    $para.Range.Font.SizeBi = 10
    if ($bold) {
        $para.Range.Font.Bold = $true
    }
}

# Anchor: paragraph 83 is "Potentially..."
$anchor = $d.Paragraphs.Item(83)

# Insert empty para (84)
$anchor.Range.InsertParagraphAfter()

# Insert "*Need better solution..." para (85)
$p84 = $d.Paragraphs.Item(84)
$p84.Range.InsertParagraphAfter()
$p85 = $d.Paragraphs.Item(85)
$p85.Range.Text = "*Need better solution for part b and c"
Set-ParaFont $p85 $false

# Insert empty para (86)
$p85.Range.InsertParagraphAfter()

# Insert "4. Evaluate Each Solution:" bold para (87)
$p86 = $d.Paragraphs.Item(86)
$p86.Range.InsertParagraphAfter()
$p87 = $d.Paragraphs.Item(87)
$p87.Range.Text = "4. Evaluate Each Solution:"
Set-ParaFont $p87 $true

# Insert "The solution for part a..." para (88)
$p87.Range.InsertParagraphAfter()
$p88 = $d.Paragraphs.Item(88)
$p88.Range.Text = "The solution for part a will meet the goal, however, it is slightly more difficult to apply to part b, and highly difficult and tedious to apply to part c."
Set-ParaFont $p88 $false

# Insert empty para (89)
$p88.Range.InsertParagraphAfter()

# Insert "*Need an evaluation..." para (90)
$p89 = $d.Paragraphs.Item(89)
$p89.Range.InsertParagraphAfter()
$p90 = $d.Paragraphs.Item(90)
$p90.Range.Text = "*Need an evaluation and better solution for part b and c"
Set-ParaFont $p90 $false

Write-Output "checkpoint B"
for ($i = 82; $i -le 96; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($pp.Range.Text)]"
}

# Handle bookmark: currently at end of paragraph 83, should move to end of paragraph 90
$bm = $d.Bookmarks.Item("_GoBack")
Write-Output "bookmark start/end before: $($bm.Start) $($bm.End)"
$p90 = $d.Paragraphs.Item(90)
Write-Output "p90 range: $($p90.Range.Start) - $($p90.Range.End)"
